# Popravila scenarij i tok dogadjaja za registriranje korisnika,
# u dijelu glavnih aktera (Korisnik -> Glavni korisnik)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Glavni tok actor-table header (row 16, column A): "Korisnik alarm sistema"
# -> "Glavni korisnik alarm sistema"
$ws.Range("A16").Value = "Glavni korisnik alarm sistema"

# Alternativni tok 1 actor-table header (row 27, column A): same change
$ws.Range("A27").Value = "Glavni korisnik alarm sistema"

# Primarni akteri (row 9, column B): "Korisnik alarm sistema i sistem za
# registraciju korisnika" -> "Glavni korisnik alarm sistema i sistem za
# registraciju korisnika"
$ws.Range("B9").Value = "Glavni korisnik alarm sistema i sistem za registraciju korisnika"

# Row 9 grew slightly taller to fit the longer actor text.
$ws.Rows(9).RowHeight = 36.75

# Update the view state to match: scrolled down a bit further, and the
# whole of row 9 (the actors row) selected instead of a single cell.
$ws.Rows(9).Select()
